$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "68.733.46"
$ws.Range("E2").Value = "  +1.14%  "
Set-TextCell $ws.Range("D3") "3.282.92"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell $ws.Range("D5") "585.21"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -0.98%  "
Set-TextCell $ws.Range("D11") "0.421"
$ws.Range("E11").Value = "  +0.54%  "
Set-TextCell $ws.Range("D12") "3.853.63"
$ws.Range("E12").Value = "  +0.16%  "
Set-TextCell $ws.Range("D14") "28.47"
$ws.Range("E14").Value = "  -0.58%  "
Set-TextCell $ws.Range("D15") "68.744.95"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("E16").Value = "  +1.36%  "
Set-TextCell $ws.Range("D17") "3.305.81"
$ws.Range("E17").Value = "  +0.91%  "
Set-TextCell $ws.Range("D18") "5.89"
$ws.Range("E18").Value = "  +0.61%  "
Set-TextCell $ws.Range("D19") "13.71"
$ws.Range("E19").Value = "  +1.08%  "
Set-TextCell $ws.Range("D20") "395.41"
$ws.Range("E20").Value = "  +4.68%  "
Set-TextCell $ws.Range("D21") "7.74"
$ws.Range("E21").Value = "  +0.57%  "
Set-TextCell $ws.Range("D22") "71.73"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  -0.02%  "
Set-TextCell $ws.Range("D24") "0.520"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +4.86%  "
Set-TextCell $ws.Range("D27") "9.83"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").Value = "  -0.06%  "
Set-TextCell $ws.Range("D29") "5.78"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E31").Value = "  +0.64%  "
Set-TextCell $ws.Range("D32") "7.18"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("E33").Value = "  +0.92%  "
Set-TextCell $ws.Range("D34") "0.999"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -1.29%  "
Set-TextCell $ws.Range("D36") "163.31"
$ws.Range("E36").Value = "  +0.62%  "
Set-TextCell $ws.Range("D37") "2.00"
$ws.Range("E37").Value = "  +8.19%  "
$ws.Range("E38").Value = "  -3.26%  "
Set-TextCell $ws.Range("D39") "26.87"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -3.36%  "
Set-TextCell $ws.Range("D43") "25.70"
$ws.Range("E43").Value = "  -0.38%  "
Set-TextCell $ws.Range("D44") "41.42"
$ws.Range("E44").Value = "  +1.34%  "
Set-TextCell $ws.Range("D45") "0.0693"
$ws.Range("E45").Value = "  +1.42%  "
Set-TextCell $ws.Range("D46") "2.662.87"
$ws.Range("E46").Value = "  -0.57%  "
Set-TextCell $ws.Range("D47") "341.49"
$ws.Range("E47").Value = "  -3.13%  "
Set-TextCell $ws.Range("D48") "0.0284"
$ws.Range("E48").Value = "  -0.25%  "
Set-TextCell $ws.Range("D49") "6.38"
$ws.Range("E49").Value = "  +3.14%  "
Set-TextCell $ws.Range("D50") "31.99"
$ws.Range("E50").Value = "  +2.52%  "
Set-TextCell $ws.Range("D51") "0.997"
$ws.Range("E51").Value = "  -0.64%  "

Write-Output "Update complete"
